# Changelog.docx - 1.7.4 - Release #2
#
# Adds two new changelog bullet paragraphs under the existing "- 1.7.4 -"
# heading:
#   1. "Update WeightsValues"  (bold)  -> inserted right before
#      "Fix tooltip messages on equipped items"
#   2. "Update LibDBIcon"      (regular) -> inserted right after
#      "Minimap icon"
#
# Both new paragraphs inherit the surrounding run formatting
# (Helvetica, sz 24 / sz-cs 24) automatically because InsertParagraphBefore/
# InsertParagraphAfter duplicate the adjoining paragraph's formatting.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, [string]$text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        $t = $p.Range.Text
        # Paragraph.Range.Text includes the trailing paragraph mark (\r)
        if ($t.TrimEnd("`r") -eq $text) {
            return $p
        }
    }
    return $null
}

# --- 1) Insert "Update WeightsValues" (bold) before "Fix tooltip messages on equipped items"
$fixTooltipPara = Find-ParagraphByText $d "Fix tooltip messages on equipped items"
$fixTooltipPara.Range.InsertParagraphBefore()

# Re-locate the paragraph by text: inserting a paragraph shifts indices, and the
# originally-held $fixTooltipPara reference keeps its stale pre-insert Index, so
# searching again (rather than trusting ".Index - 1") is what actually lands on
# the freshly created empty paragraph.
$fixTooltipPara = Find-ParagraphByText $d "Fix tooltip messages on equipped items"
$newPara1 = $d.Paragraphs($fixTooltipPara.Index - 1)
$newRange1 = $d.Range($newPara1.Range.Start, $newPara1.Range.End - 1)
$newRange1.Text = "Update WeightsValues"
$newRange1.Font.Bold = 1

# --- 2) Insert "Update LibDBIcon" (not bold) after "Minimap icon"
$minimapPara = Find-ParagraphByText $d "Minimap icon"
$minimapPara.Range.InsertParagraphAfter()

$minimapPara = Find-ParagraphByText $d "Minimap icon"
$newPara2 = $d.Paragraphs($minimapPara.Index + 1)
$newRange2 = $d.Range($newPara2.Range.Start, $newPara2.Range.End - 1)
$newRange2.Text = "Update LibDBIcon"
